$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells in row 1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting from the last existing header cell (AC1)
# so the new header cells share the same bold/centered/bordered style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the Wins/Losses/Ties team-record columns for every data row (2-49)
$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 76  # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 85  # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF -> Ties
}
